# Update the cryptocurrency Price (D) and Volume(1h) (E) columns with refreshed values.
# A leading apostrophe is used for D-column values that look like plain numbers so Excel
# keeps them stored as text (matching the source data, e.g. "1.00" rather than 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.070.34'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '3.538.62'
$ws.Range("E3").Value = '  +4.36%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''598.52'
$ws.Range("E5").Value = '  +3.61%  '
$ws.Range("D6").Value = '''137.53'
$ws.Range("D7").Value = '3.538.21'
$ws.Range("E7").Value = '  +4.35%  '
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("D9").Value = '''0.493'
$ws.Range("E9").Value = '  +0.79%  '
$ws.Range("E10").Value = '  +3.87%  '
$ws.Range("E11").Value = '  -0.98%  '
$ws.Range("D12").Value = '''0.386'
$ws.Range("E12").Value = '  +4.72%  '
$ws.Range("D13").Value = '4.143.36'
$ws.Range("E13").Value = '  +4.41%  '
$ws.Range("E14").Value = '  +4.64%  '
$ws.Range("D15").Value = '''27.30'
$ws.Range("E15").Value = '  +5.81%  '
$ws.Range("D16").Value = '3.537.08'
$ws.Range("E16").Value = '  +4.12%  '
$ws.Range("E17").Value = '  +1.68%  '
$ws.Range("D18").Value = '65.091.30'
$ws.Range("E18").Value = '  +0.45%  '
$ws.Range("D19").Value = '''10.13'
$ws.Range("E19").Value = '  +7.18%  '
$ws.Range("D20").Value = '''5.87'
$ws.Range("E20").Value = '  +2.46%  '
$ws.Range("D21").Value = '''14.19'
$ws.Range("E21").Value = '  +6.77%  '
$ws.Range("D22").Value = '''390.08'
$ws.Range("E22").Value = '  +3.73%  '
$ws.Range("D23").Value = '''0.574'
$ws.Range("E23").Value = '  +5.77%  '
$ws.Range("D24").Value = '3.682.43'
$ws.Range("E24").Value = '  +4.25%  '
$ws.Range("D25").Value = '''73.89'
$ws.Range("E25").Value = '  +3.36%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").Value = '''0.0000112'
$ws.Range("E27").Value = '  +9.98%  '
$ws.Range("D28").Value = '''7.83'
$ws.Range("E28").Value = '  +13.50%  '
$ws.Range("D29").Value = '''0.983'
$ws.Range("E29").Value = '  -1.93%  '
$ws.Range("D30").Value = '''2.27'
$ws.Range("E30").Value = '  +4.50%  '
$ws.Range("D31").Value = '''8.27'
$ws.Range("E31").Value = '  +5.00%  '
$ws.Range("D32").Value = '3.557.40'
$ws.Range("E32").Value = '  +4.56%  '
$ws.Range("E33").Value = '  +26.11%  '
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("D35").Value = '''23.87'
$ws.Range("E35").Value = '  +4.73%  '
$ws.Range("D36").Value = '''0.144'
$ws.Range("E36").Value = '  +3.21%  '
$ws.Range("E37").Value = '  +11.33%  '
$ws.Range("D38").Value = '''6.94'
$ws.Range("E38").Value = '  +5.49%  '
$ws.Range("D39").Value = '''168.47'
$ws.Range("E39").Value = '  -0.71%  '
$ws.Range("D40").Value = '''5.02'
$ws.Range("E40").Value = '  +10.13%  '
$ws.Range("D41").Value = '''0.0805'
$ws.Range("E41").Value = '  +8.52%  '
$ws.Range("D42").Value = '''0.823'
$ws.Range("E42").Value = '  +2.48%  '
$ws.Range("D43").Value = '''26.59'
$ws.Range("E43").Value = '  +22.97%  '
$ws.Range("D44").Value = '''1.00'
$ws.Range("E44").Value = '  -0.10%  '
$ws.Range("D45").Value = '''42.46'
$ws.Range("E45").Value = '  +0.45%  '
$ws.Range("D46").Value = '''4.44'
$ws.Range("E46").Value = '  +4.15%  '
$ws.Range("E47").Value = '  +12.03%  '
$ws.Range("E48").Value = '  +6.99%  '
$ws.Range("D49").Value = '''6.84'
$ws.Range("E49").Value = '  +7.18%  '
$ws.Range("D50").Value = '2.399.17'
$ws.Range("E50").Value = '  +12.03%  '
$ws.Range("D51").Value = '''309.29'
$ws.Range("E51").Value = '  +18.53%  '
